$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D15").Value = "2016-03-04 15:41:41"
$wsZhCn.Range("G15").Value = "2016-03-04 15:42:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D15").Value = "2016-03-04 15:41:53"
$wsDeDe.Range("G15").Value = "2016-03-04 15:42:49"
